# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7 updates
$ws.Range("G7").Value = 1.7
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.4
$ws.Range("L7").Value = 5.5
$ws.Range("N7").Value = 7.5
$ws.Range("Z7").Value = 13
$ws.Range("AI7").Value = 23
$ws.Range("AJ7").Value = 17
$ws.Range("AO7").Value = 9.5
$ws.Range("AR7").Value = 67
$ws.Range("AW7").Value = 6.5

# Row 8 updates
$ws.Range("N8").Value = 7.9
